# Resize the inline picture ("Picture") to its native 96-DPI size
# (5in x 4in == 360pt x 288pt == 4572000 x 3657600 EMU) and lock its
# aspect ratio so subsequent resizes preserve proportions.

$d = $word.ActiveDocument

$shape = $d.InlineShapes(1)
$shape.LockAspectRatio = 1
$shape.Width = 360
$shape.Height = 288
